$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$questionsText = 'questions = [
    {
        "title": "You ran a MATLAB script that created multiple variables. You now want to remove all variables from the workspace to free up memory. Which of these statements can you use to accomplish this?",
        "ques_type": 15,
        "options": [
            "clc",
            "clear",
            "clear all",
            "delete",
            "clear vars"
        ],
        "score": [
            "clear",
            "clear all"
        ]
    },
    {
        "title": "Consider the following equation: You want to convert this equation into MATLAB code. Which of these corresponds to the correct translation of the equation?",
        "ques_type": 2,
        "options": [
            "a/b+c^2+c/d*e/f",
            "a/(b+c^2)+c/d*e/f",
            "a/(b+c*2)+c/d*e/f",
            "a/b+c**2+c/d*e/f"
        ],
        "score": "a/(b+c^2)+c/d*e/f"
    },
    {
        "title": "You are provided with the following piece of code: temp = 4 \nwhile (temp &lt 9) \ndisp(''Hello World!'') \ntemp = temp + 2.5 \nend\n In the given code, how many times will the display command be executed?",
        "ques_type": 2,
        "options": [
            "2",
            "3",
            "4",
            "9"
        ],
        "score": "2"
    },
    {
        "title": "You are provided with two vectors: \u201cage\u201d and \u201cweight.\u201d You want to create a line plot, where age corresponds to the x-axis and weight corresponds to the y-axis.  Which of these statements should you use to accomplish this?",
        "ques_type": 2,
        "options": [
            "plot(age,weight)",
            "plot(weight,age)",
            "plot[(age,weight)]",
            "plot([weight,age])"
        ],
        "score": "plot(age,weight)"
    }
]'

$ws.Rows("2:2").Delete()
$ws.Range("A1").Value = $questionsText
